# Atualização de bases das ligas, do dia: 28-05-2024 às 20:56
# Swap the data (columns B:AD, i.e. everything except the row-index
# column A) between each of the following row pairs. This corrects rows
# whose underlying match records had been placed on the wrong line while
# the "id" ordinal in column A stays tied to its original row position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(22, 23),
    @(41, 42),
    @(58, 59),
    @(91, 92),
    @(135, 136)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1" + ":AD$r1")
    $range2 = $ws.Range("B$r2" + ":AD$r2")

    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value2 = $values2
    $range2.Value2 = $values1
}
